# Update anonymized contract outputs: renumber the PERSON_## placeholders
# inside the skloňování (declension) example list. Each bullet paragraph is a
# single run; we rewrite each paragraph's visible text in place (by its
# Paragraphs() index) so the paragraph mark / numbering / formatting stay
# untouched. Using positional indices (rather than Find/Replace on the old
# text) avoids any ambiguity from the fact that several "before" strings
# equal other entries' "after" strings once the shift has been applied.

$d = $word.ActiveDocument

$replacements = @(
    @{ Index = 39; Old = "[[PERSON_41]] – „pro [[PERSON_41]]“, „o [[PERSON_42]]“"; New = "[[PERSON_41]] – „pro [[PERSON_42]]“, „o [[PERSON_43]]“" },
    @{ Index = 40; Old = "[[PERSON_43]] – „k [[PERSON_43]]“, „o [[PERSON_43]]“";   New = "[[PERSON_44]] – „k [[PERSON_44]]“, „o [[PERSON_44]]“" },
    @{ Index = 41; Old = "[[PERSON_44]] – „o [[PERSON_44]]“, „s [[PERSON_44]]“";  New = "[[PERSON_45]] – „o [[PERSON_45]]“, „s [[PERSON_45]]“" },
    @{ Index = 42; Old = "[[PERSON_45]] – „s [[PERSON_45]]“, „o [[PERSON_45]]“";  New = "[[PERSON_46]] – „s [[PERSON_46]]“, „o [[PERSON_46]]“" },
    @{ Index = 43; Old = "[[PERSON_46]] – „s [[PERSON_46]]“, „o [[PERSON_46]]“";  New = "[[PERSON_47]] – „s [[PERSON_47]]“, „o [[PERSON_47]]“" },
    @{ Index = 44; Old = "[[PERSON_47]] – „u [[PERSON_48]]“, „o [[PERSON_47]]“";  New = "[[PERSON_48]] – „u [[PERSON_48]]“, „o [[PERSON_48]]“" },
    @{ Index = 46; Old = "[[PERSON_50]] – „o [[PERSON_51]]“, „s [[PERSON_50]]“";  New = "[[PERSON_50]] – „o [[PERSON_50]]“, „s [[PERSON_50]]“" },
    @{ Index = 47; Old = "[[PERSON_52]] – „k [[PERSON_52]]“, „o [[PERSON_52]]“";  New = "[[PERSON_51]] – „k [[PERSON_51]]“, „o [[PERSON_51]]“" },
    @{ Index = 48; Old = "[[PERSON_53]] – „o [[PERSON_53]]“, „s [[PERSON_53]]“";  New = "[[PERSON_52]] – „o [[PERSON_52]]“, „s [[PERSON_52]]“" },
    @{ Index = 49; Old = "[[PERSON_54]] – „s [[PERSON_54]]“, „o [[PERSON_54]]“";  New = "[[PERSON_53]] – „s [[PERSON_53]]“, „o [[PERSON_53]]“" },
    @{ Index = 50; Old = "[[PERSON_55]] – „s [[PERSON_55]]“, „o [[PERSON_55]]“";  New = "[[PERSON_54]] – „s [[PERSON_54]]“, „o [[PERSON_54]]“" },
    @{ Index = 51; Old = "[[PERSON_56]] – „o [[PERSON_56]]“, „s [[PERSON_56]]“";  New = "[[PERSON_55]] – „o [[PERSON_55]]“, „s [[PERSON_55]]“" },
    @{ Index = 52; Old = "[[PERSON_57]] – „s [[PERSON_57]]“, „o [[PERSON_57]]“";  New = "[[PERSON_56]] – „s [[PERSON_56]]“, „o [[PERSON_56]]“" },
    @{ Index = 53; Old = "[[PERSON_58]] – „o [[PERSON_58]]“, „s [[PERSON_58]]“";  New = "[[PERSON_57]] – „o [[PERSON_57]]“, „s [[PERSON_57]]“" },
    @{ Index = 54; Old = "[[PERSON_59]] – „s [[PERSON_59]]“, „o [[PERSON_59]]“";  New = "[[PERSON_58]] – „s [[PERSON_58]]“, „o [[PERSON_58]]“" }
)

foreach ($rep in $replacements) {
    $para = $d.Paragraphs($rep.Index)
    $current = $para.Range.Text
    # Trim trailing paragraph-mark / cell-mark characters (CR = 0x0D) for comparison.
    $currentTrimmed = $current.TrimEnd([char]13, [char]7)
    if ($currentTrimmed -ne $rep.Old) {
        # Fallback: the expected paragraph wasn't at this index (e.g. layout
        # shifted) - locate it by its exact original text anywhere in the
        # document instead, so the edit still lands correctly.
        $search = $d.Content
        $ok = $search.Find.Execute($rep.Old, $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
        if ($ok) {
            $search.Text = $rep.New
        } else {
            Write-Output "WARNING: could not locate paragraph $($rep.Index) / text [$($rep.Old)] (saw [$currentTrimmed])"
        }
    } else {
        $para.Range.Text = $rep.New
    }
}

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
